{"js": "const body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\n// Build a lookup of paragraphs keyed by their exact current text.\n// (Each lookup uses the ORIGINAL pre-edit text, so this is computed once\n// up front before any mutations take place.)\nconst byText = {};\nfor (const p of paras.items) {\n  byText[p.text] = p;\n}\n\n// 1) Merge the two \"Direction: \" / \"Dream Sequence\" runs into a single run.\nconst directionPara = byText[\"Direction: Dream Sequence\"];\nif (directionPara) {\n  directionPara.insertText(\"Direction: Dream Sequence\", Word.InsertLocation.replace);\n}\n\n// 2) Split the \"My heart rate slows down...\" paragraph into three paragraphs.\nconst heartRatePara = byText[\n  \"My heart rate slows down as I fumble around for my phone to check whether or not I can sneak a few extra minutes of sleep in. Unfortunately, it turns out that if I don\\u2019t get up now, I\\u2019ll most likely be late for school.\"\n];\nif (heartRatePara) {\n  heartRatePara.insertText(\n    \"My heart rate slows down as I fumble around for my phone, wanting to check whether or not I can sneak a few extra minutes of sleep in.\",\n    Word.InsertLocation.replace\n  );\n  heartRatePara.insertParagraph(\n    \"Ah. If I don\\u2019t get up now, I\\u2019ll most likely be late for school. Unfortunate.\",\n    Word.InsertLocation.after\n  );\n}\n\nconst fallAsleepPara = byText[\n  \"Ah, well. I guess I probably wouldn\\u2019t have been able to fall asleep again anyways.\"\n];\nif (fallAsleepPara) {\n  fallAsleepPara.insertText(\n    \"I guess I probably wouldn\\u2019t have been able to fall asleep again anyways. Not after that.\",\n    Word.InsertLocation.replace\n  );\n}\n\n// 3) Mom's greeting line: \"waving smile\" -> \"neutral smiling\".\nconst momGreetingPara = byText[\n  \"Mom (waving smile): Oh, I was just about to go and wake you up. Good morning.\"\n];\nif (momGreetingPara) {\n  momGreetingPara.insertText(\n    \"Mom (neutral smiling): Oh, I was just about to go and wake you up. Good morning.\",\n    Word.InsertLocation.replace\n  );\n}\n\n// 4) Insert a new \"Mom (neutral neutral):\" paragraph before \"I sit down at the table...\"\nconst sitDownPara = byText[\"I sit down at the table, where, as always, breakfast waits for me.\"];\nif (sitDownPara) {\n  sitDownPara.insertParagraph(\"Mom (neutral neutral):\", Word.InsertLocation.before);\n}\n\n// 5) \"smiling_nervous\" -> \"worried_slightly\" for the \"You look tired\" line.\nconst lookTiredPara = byText[\"Mom (neutral smiling_nervous): You look tired. Did you sleep alright?\"];\nif (lookTiredPara) {\n  lookTiredPara.insertText(\n    \"Mom (neutral worried_slightly): You look tired. Did you sleep alright?\",\n    Word.InsertLocation.replace\n  );\n}\n\n// 6) \"smiling_nervous\" -> \"thinking\" for the \"You went to bed really early\" line,\n//    then insert a new \"Mom (neutral worried):\" paragraph right after it.\nconst wentToBedPara = byText[\"Mom (neutral smiling_nervous): You went to bed really early last night, too.\"];\nif (wentToBedPara) {\n  wentToBedPara.insertText(\n    \"Mom (neutral thinking): You went to bed really early last night, too.\",\n    Word.InsertLocation.replace\n  );\n  wentToBedPara.insertParagraph(\"Mom (neutral worried):\", Word.InsertLocation.after);\n}\n\n// 7) \"Mom (neutral worried): Don't push yourself too hard, okay?\" -> \"smiling_worried\".\nconst pushYourselfPara = byText[\"Mom (neutral worried): Don\\u2019t push yourself too hard, okay?\"];\nif (pushYourselfPara) {\n  pushYourselfPara.insertText(\n    \"Mom (neutral smiling_worried): Don\\u2019t push yourself too hard, okay?\",\n    Word.InsertLocation.replace\n  );\n}\n\n// 8) \"Mom (neutral sigh): When it comes to school...\" -> \"smiling_nervous\".\nconst schoolPara = byText[\n  \"Mom (neutral sigh): When it comes to school you could probably push yourself a bit harder.\"\n];\nif (schoolPara) {\n  schoolPara.insertText(\n    \"Mom (neutral smiling_nervous): When it comes to school you could probably push yourself a bit harder.\",\n    Word.InsertLocation.replace\n  );\n}\n\n// 9) Insert a new \"Mom (neutral smiling):\" paragraph before \"Pro: Uh...\" (after the\n//    \"Who's gonna support me when I retire?\" line).\nconst whoGonnaSupportPara = byText[\"Mom (neutral sigh): Who\\u2019s gonna support me when I retire?\"];\nif (whoGonnaSupportPara) {\n  whoGonnaSupportPara.insertParagraph(\"Mom (neutral smiling):\", Word.InsertLocation.after);\n}\n\n// 10) \"Mom (neutral smiling): I'm just kidding...\" -> \"smiling_eyes_closed\".\nconst justKiddingPara = byText[\"Mom (neutral smiling): I\\u2019m just kidding. Just do your best.\"];\nif (justKiddingPara) {\n  justKiddingPara.insertText(\n    \"Mom (neutral smiling_eyes_closed): I\\u2019m just kidding. Just do your best.\",\n    Word.InsertLocation.replace\n  );\n}\n\n// 11) Remove the trailing space on the final line of the document.\nconst livelyGirlPara = byText[\"And with that, she spins around and trots off. What a lively girl. \"];\nif (livelyGirlPara) {\n  livelyGirlPara.insertText(\n    \"And with that, she spins around and trots off. What a lively girl.\",\n    Word.InsertLocation.replace\n  );\n}\n\nawait context.sync();\n", "ps1": "function Get-ParaByText($doc, $text) {\n    foreach ($p in $doc.Paragraphs) {\n        if ($p.Range.Text.TrimEnd(\"`r\") -eq $text) {\n            return $p\n        }\n    }\n    return $null\n}\n\nfunction Replace-ParaText($doc, $oldText, $newText) {\n    $find = $doc.Content.Find\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $result = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n\n$d = $word.ActiveDocument\n\n# 1) Merge the \"Direction: \" / \"Dream Sequence\" runs into a single run.\nReplace-ParaText $d \"Direction: Dream Sequence\" \"Direction: Dream Sequence\"\n\n# 2) Split the \"My heart rate slows down...\" paragraph into three paragraphs.\nReplace-ParaText $d \"My heart rate slows down as I fumble around for my phone to check whether or not I can sneak a few extra minutes of sleep in. Unfortunately, it turns out that if I don\u2019t get up now, I\u2019ll most likely be late for school.\" \"My heart rate slows down as I fumble around for my phone, wanting to check whether or not I can sneak a few extra minutes of sleep in.\"\n\n$heartRatePara = Get-ParaByText $d \"My heart rate slows down as I fumble around for my phone, wanting to check whether or not I can sneak a few extra minutes of sleep in.\"\nif ($heartRatePara -ne $null) {\n    $heartRatePara.Range.InsertParagraphAfter()\n    $newPara = $heartRatePara.Next()\n    $newPara.Range.Text = \"Ah. If I don\u2019t get up now, I\u2019ll most likely be late for school. Unfortunate.\"\n}\n\nReplace-ParaText $d \"Ah, well. I guess I probably wouldn\u2019t have been able to fall asleep again anyways.\" \"I guess I probably wouldn\u2019t have been able to fall asleep again anyways. Not after that.\"\n\n# 3) Mom's greeting line: \"waving smile\" -> \"neutral smiling\".\nReplace-ParaText $d \"Mom (waving smile): Oh, I was just about to go and wake you up. Good morning.\" \"Mom (neutral smiling): Oh, I was just about to go and wake you up. Good morning.\"\n\n# 4) Insert a new \"Mom (neutral neutral):\" paragraph before \"I sit down at the table...\"\n# NOTE: After InsertParagraphBefore(), the paragraph object itself ends up referring to\n# the newly-created (now preceding) empty paragraph, so we set its text directly.\n$sitDownPara = Get-ParaByText $d \"I sit down at the table, where, as always, breakfast waits for me.\"\nif ($sitDownPara -ne $null) {\n    $sitDownPara.Range.InsertParagraphBefore()\n    $sitDownPara.Range.Text = \"Mom (neutral neutral):\"\n}\n\n# 5) \"smiling_nervous\" -> \"worried_slightly\" for the \"You look tired\" line.\nReplace-ParaText $d \"Mom (neutral smiling_nervous): You look tired. Did you sleep alright?\" \"Mom (neutral worried_slightly): You look tired. Did you sleep alright?\"\n\n# 6) \"smiling_nervous\" -> \"thinking\" for the \"You went to bed really early\" line,\n#    then insert a new \"Mom (neutral worried):\" paragraph right after it.\nReplace-ParaText $d \"Mom (neutral smiling_nervous): You went to bed really early last night, too.\" \"Mom (neutral thinking): You went to bed really early last night, too.\"\n\n$wentToBedPara = Get-ParaByText $d \"Mom (neutral thinking): You went to bed really early last night, too.\"\nif ($wentToBedPara -ne $null) {\n    $wentToBedPara.Range.InsertParagraphAfter()\n    $newPara = $wentToBedPara.Next()\n    $newPara.Range.Text = \"Mom (neutral worried):\"\n}\n\n# 7) \"Mom (neutral worried): Don't push yourself too hard, okay?\" -> \"smiling_worried\".\nReplace-ParaText $d \"Mom (neutral worried): Don\u2019t push yourself too hard, okay?\" \"Mom (neutral smiling_worried): Don\u2019t push yourself too hard, okay?\"\n\n# 8) \"Mom (neutral sigh): When it comes to school...\" -> \"smiling_nervous\".\nReplace-ParaText $d \"Mom (neutral sigh): When it comes to school you could probably push yourself a bit harder.\" \"Mom (neutral smiling_nervous): When it comes to school you could probably push yourself a bit harder.\"\n\n# 9) Insert a new \"Mom (neutral smiling):\" paragraph after \"Who's gonna support me when I retire?\"\n$whoGonnaSupportPara = Get-ParaByText $d \"Mom (neutral sigh): Who\u2019s gonna support me when I retire?\"\nif ($whoGonnaSupportPara -ne $null) {\n    $whoGonnaSupportPara.Range.InsertParagraphAfter()\n    $newPara = $whoGonnaSupportPara.Next()\n    $newPara.Range.Text = \"Mom (neutral smiling):\"\n}\n\n# 10) \"Mom (neutral smiling): I'm just kidding...\" -> \"smiling_eyes_closed\".\nReplace-ParaText $d \"Mom (neutral smiling): I\u2019m just kidding. Just do your best.\" \"Mom (neutral smiling_eyes_closed): I\u2019m just kidding. Just do your best.\"\n\n# 11) Remove the trailing space on the final line of the document.\nReplace-ParaText $d \"And with that, she spins around and trots off. What a lively girl. \" \"And with that, she spins around and trots off. What a lively girl.\"\n"}
